$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels for columns Q and R (row 1)
$ws.Range("Q1").Value = "Yglm"
$ws.Range("R1").Value = "Yerr"

# Update recomputed values for columns Q-V across data rows 2-33
# Row 2
$ws.Range("Q2").Value = 0.6876816604375845
$ws.Range("R2").Value = -0.03675689776427005
$ws.Range("S2").Value = [double]"1.285869558979395e-13"
$ws.Range("T2").Value = 0.6876816604377131
$ws.Range("U2").Value = 0.001834048387843977
$ws.Range("V2").Value = 0.0005111754801311683
# Row 3
$ws.Range("Q3").Value = 0.485028243619635
$ws.Range("R3").Value = 0.1713731834993147
$ws.Range("S3").Value = [double]"1.285869558979395e-13"
$ws.Range("T3").Value = 0.4850282436197636
$ws.Range("U3").Value = 0.001834048387843977
$ws.Range("V3").Value = 0.0005111754801311683
# Row 4
$ws.Range("Q4").Value = 0.4380124470085797
$ws.Range("R4").Value = -0.07592273045324749
$ws.Range("S4").Value = [double]"1.285869558979395e-13"
$ws.Range("T4").Value = 0.4380124470087082
$ws.Range("U4").Value = 0.001834048387843977
$ws.Range("V4").Value = 0.0005111754801311683
# Row 5
$ws.Range("Q5").Value = 0.4666791532827632
$ws.Range("R5").Value = -0.111862556459848
$ws.Range("S5").Value = [double]"1.285869558979395e-13"
$ws.Range("T5").Value = 0.4666791532828918
$ws.Range("U5").Value = 0.001834048387843977
$ws.Range("V5").Value = 0.0005111754801311683
# Row 6
$ws.Range("Q6").Value = 0.3835093432858883
$ws.Range("R6").Value = -0.04259233072114649
$ws.Range("S6").Value = [double]"1.285869558979395e-13"
$ws.Range("T6").Value = 0.3835093432860169
$ws.Range("U6").Value = 0.001834048387843977
$ws.Range("V6").Value = 0.0005111754801311683
# Row 7
$ws.Range("Q7").Value = 0.361573517919367
$ws.Range("R7").Value = -0.0150001482877512
$ws.Range("S7").Value = [double]"1.285869558979395e-13"
$ws.Range("T7").Value = 0.3615735179194955
$ws.Range("U7").Value = 0.001834048387843977
$ws.Range("V7").Value = 0.0005111754801311683
# Row 8
$ws.Range("Q8").Value = 0.3328663257406512
$ws.Range("R8").Value = 0.01361923288485356
$ws.Range("S8").Value = [double]"1.285869558979395e-13"
$ws.Range("T8").Value = 0.3328663257407797
$ws.Range("U8").Value = 0.001834048387843977
$ws.Range("V8").Value = 0.0005111754801311683
# Row 9
$ws.Range("Q9").Value = 0.2735951854887018
$ws.Range("R9").Value = 0.07284445240701665
$ws.Range("S9").Value = [double]"1.285869558979395e-13"
$ws.Range("T9").Value = 0.2735951854888303
$ws.Range("U9").Value = 0.001834048387843977
$ws.Range("V9").Value = 0.0005111754801311683
# Row 10
$ws.Range("Q10").Value = 0.3458046741187144
$ws.Range("R10").Value = 0.007776640788102973
$ws.Range("S10").Value = [double]"1.285869558979395e-13"
$ws.Range("T10").Value = 0.345804674118843
$ws.Range("U10").Value = 0.001834048387843977
$ws.Range("V10").Value = 0.0005111754801311683
# Row 11
$ws.Range("Q11").Value = 0.3062128923254483
$ws.Range("R11").Value = 0.04386926338254749
$ws.Range("S11").Value = [double]"1.285869558979395e-13"
$ws.Range("T11").Value = 0.3062128923255769
$ws.Range("U11").Value = 0.001834048387843977
$ws.Range("V11").Value = 0.0005111754801311683
# Row 12
$ws.Range("Q12").Value = 0.3598323936330917
$ws.Range("R12").Value = 0.001742530581442137
$ws.Range("S12").Value = [double]"1.285869558979395e-13"
$ws.Range("T12").Value = 0.3598323936332202
$ws.Range("U12").Value = 0.001834048387843977
$ws.Range("V12").Value = 0.0005111754801311683
# Row 13
$ws.Range("Q13").Value = 0.6271121682286698
$ws.Range("R13").Value = 0.02563479707444238
$ws.Range("S13").Value = [double]"1.285869558979395e-13"
$ws.Range("T13").Value = 0.6271121682287983
$ws.Range("U13").Value = 0.001834048387843977
$ws.Range("V13").Value = 0.0005111754801311683
# Row 14
$ws.Range("Q14").Value = 0.4153933263559565
$ws.Range("R14").Value = -0.062204183088575
$ws.Range("S14").Value = [double]"1.285869558979395e-13"
$ws.Range("T14").Value = 0.4153933263560851
$ws.Range("U14").Value = 0.001834048387843977
$ws.Range("V14").Value = 0.0005111754801311683
# Row 15
$ws.Range("Q15").Value = 0.6197720490325688
$ws.Range("R15").Value = 0.03511488132287299
$ws.Range("S15").Value = [double]"1.285869558979395e-13"
$ws.Range("T15").Value = 0.6197720490326973
$ws.Range("U15").Value = 0.001834048387843977
$ws.Range("V15").Value = 0.0005111754801311683
# Row 16
$ws.Range("Q16").Value = 0.3950080504618158
$ws.Range("R16").Value = -0.03120957534001023
$ws.Range("S16").Value = [double]"1.285869558979395e-13"
$ws.Range("T16").Value = 0.3950080504619443
$ws.Range("U16").Value = 0.001834048387843977
$ws.Range("V16").Value = 0.0005111754801311683
# Row 17
$ws.Range("Q17").Value = 0.5065516705575808
$ws.Range("R17").Value = 0.09966784314690191
$ws.Range("S17").Value = [double]"1.285869558979395e-13"
$ws.Range("T17").Value = 0.5065516705577093
$ws.Range("U17").Value = 0.001834048387843977
$ws.Range("V17").Value = 0.0005111754801311683
# Row 18
$ws.Range("Q18").Value = 0.5300108218750932
$ws.Range("R18").Value = 0.1067721291424536
$ws.Range("S18").Value = [double]"1.285869558979395e-13"
$ws.Range("T18").Value = 0.5300108218752217
$ws.Range("U18").Value = 0.001834048387843977
$ws.Range("V18").Value = 0.0005111754801311683
# Row 19
$ws.Range("Q19").Value = 0.3162730747708147
$ws.Range("R19").Value = 0.03569490663255165
$ws.Range("S19").Value = [double]"1.285869558979395e-13"
$ws.Range("T19").Value = 0.3162730747709432
$ws.Range("U19").Value = 0.001834048387843977
$ws.Range("V19").Value = 0.0005111754801311683
# Row 20
$ws.Range("Q20").Value = 0.6141401951586175
$ws.Range("R20").Value = -0.01143930695691719
$ws.Range("S20").Value = [double]"1.285869558979395e-13"
$ws.Range("T20").Value = 0.6141401951587461
$ws.Range("U20").Value = 0.001834048387843977
$ws.Range("V20").Value = 0.0005111754801311683
# Row 21
$ws.Range("Q21").Value = 0.443574015457025
$ws.Range("R21").Value = -0.04156150526721536
$ws.Range("S21").Value = [double]"1.285869558979395e-13"
$ws.Range("T21").Value = 0.4435740154571535
$ws.Range("U21").Value = 0.001834048387843977
$ws.Range("V21").Value = 0.0005111754801311683
# Row 22
$ws.Range("Q22").Value = 0.4577184545255421
$ws.Range("R22").Value = [double]"-5.391559499212661e-05"
$ws.Range("S22").Value = [double]"1.285869558979395e-13"
$ws.Range("T22").Value = 0.4577184545256707
$ws.Range("U22").Value = 0.001834048387843977
$ws.Range("V22").Value = 0.0005111754801311683
# Row 23
$ws.Range("Q23").Value = 0.3718069884886606
$ws.Range("R23").Value = -0.02186850880979829
$ws.Range("S23").Value = [double]"1.285869558979395e-13"
$ws.Range("T23").Value = 0.3718069884887892
$ws.Range("U23").Value = 0.001834048387843977
$ws.Range("V23").Value = 0.0005111754801311683
# Row 24
$ws.Range("Q24").Value = 0.7492466418763813
$ws.Range("R24").Value = -0.1008053388724446
$ws.Range("S24").Value = [double]"1.285869558979395e-13"
$ws.Range("T24").Value = 0.7492466418765099
$ws.Range("U24").Value = 0.001834048387843977
$ws.Range("V24").Value = 0.0005111754801311683
# Row 25
$ws.Range("Q25").Value = 0.7444824352314566
$ws.Range("R25").Value = -0.09291436240159312
$ws.Range("S25").Value = [double]"1.285869558979395e-13"
$ws.Range("T25").Value = 0.7444824352315852
$ws.Range("U25").Value = 0.001834048387843977
$ws.Range("V25").Value = 0.0005111754801311683
# Row 26
$ws.Range("Q26").Value = 0.4554077105400403
$ws.Range("R26").Value = 0.1327646750056338
$ws.Range("S26").Value = [double]"1.285869558979395e-13"
$ws.Range("T26").Value = 0.4554077105401689
$ws.Range("U26").Value = 0.001834048387843977
$ws.Range("V26").Value = 0.0005111754801311683
# Row 27
$ws.Range("Q27").Value = 0.6588836812540453
$ws.Range("R27").Value = -0.02119904567370978
$ws.Range("S27").Value = [double]"1.285869558979395e-13"
$ws.Range("T27").Value = 0.6588836812541738
$ws.Range("U27").Value = 0.001834048387843977
$ws.Range("V27").Value = 0.0005111754801311683
# Row 28
$ws.Range("Q28").Value = 0.6314449715686675
$ws.Range("R28").Value = 0.01430848998996692
$ws.Range("S28").Value = [double]"1.285869558979395e-13"
$ws.Range("T28").Value = 0.6314449715687961
$ws.Range("U28").Value = 0.001834048387843977
$ws.Range("V28").Value = 0.0005111754801311683
# Row 29
$ws.Range("Q29").Value = 0.5617749785344494
$ws.Range("R29").Value = 0.08301012488674875
$ws.Range("S29").Value = [double]"1.285869558979395e-13"
$ws.Range("T29").Value = 0.5617749785345779
$ws.Range("U29").Value = 0.001834048387843977
$ws.Range("V29").Value = 0.0005111754801311683
# Row 30
$ws.Range("Q30").Value = 0.3157116145950068
$ws.Range("R30").Value = 0.09542441923041139
$ws.Range("S30").Value = [double]"1.285869558979395e-13"
$ws.Range("T30").Value = 0.3157116145951354
$ws.Range("U30").Value = 0.001834048387843977
$ws.Range("V30").Value = 0.0005111754801311683
# Row 31
$ws.Range("Q31").Value = 0.7079069580402964
$ws.Range("R31").Value = -0.07474172902661091
$ws.Range("S31").Value = [double]"1.285869558979395e-13"
$ws.Range("T31").Value = 0.7079069580404249
$ws.Range("U31").Value = 0.001834048387843977
$ws.Range("V31").Value = 0.0005111754801311683
# Row 32
$ws.Range("Q32").Value = 0.3121422412936599
$ws.Range("R32").Value = 0.04082321216197832
$ws.Range("S32").Value = [double]"1.285869558979395e-13"
$ws.Range("T32").Value = 0.3121422412937884
$ws.Range("U32").Value = 0.001834048387843977
$ws.Range("V32").Value = 0.0005111754801311683
# Row 33
$ws.Range("Q33").Value = 0.4811773804434596
$ws.Range("R33").Value = 0.04777381443137646
$ws.Range("S33").Value = [double]"1.285869558979395e-13"
$ws.Range("T33").Value = 0.4811773804435882
$ws.Range("U33").Value = 0.001834048387843977
$ws.Range("V33").Value = 0.0005111754801311683
